# Auto-generated Excel COM-interop script
# Updates the cryptos worksheet with the new price/volume/ordering data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text format
# first, otherwise Excel would auto-convert the inline string into a numeric value.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.552.29"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").Value = "1.830.22"
$ws.Range("E3").Value = "  +5.17%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "344.59"
$ws.Range("E5").Value = "  +3.50%  "
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.3835"
$ws.Range("E7").Value = "  +2.49%  "
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "50.43"
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.3549"
$ws.Range("E9").Value = "  +4.20%  "
$ws.Range("D10").Value = "1.239"
$ws.Range("E10").Value = "  +3.47%  "
$ws.Range("D11").Value = "0.07797"
$ws.Range("E11").Value = "  +4.34%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "22.38"
$ws.Range("E13").Value = "  +9.97%  "
$ws.Range("D14").Value = "6.625"
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").Value = "1.826.23"
$ws.Range("E15").Value = "  +5.07%  "
$ws.Range("D16").Value = "7.252"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").Value = "0.00001130"
$ws.Range("E17").Value = "  +3.90%  "
$ws.Range("D18").Value = "0.06734"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Value = "86.93"
$ws.Range("E19").Value = "  +4.05%  "
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "17.65"
$ws.Range("E21").Value = "  +5.68%  "
$ws.Range("D22").Value = "6.584"
$ws.Range("E22").Value = "  +6.43%  "
$ws.Range("D23").Value = "13.21"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").Value = "27.574.56"
$ws.Range("E24").Value = "  +3.38%  "
$ws.Range("D25").Value = "2.471"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").Value = "2.756"
$ws.Range("E26").Value = "  +11.34%  "
$ws.Range("D27").Value = "22.23"
$ws.Range("E27").Value = "  +13.27%  "
$ws.Range("D28").Value = "1.476"
$ws.Range("E28").Value = "  +4.47%  "
$ws.Range("D29").Value = "153.63"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "2.035.23"
$ws.Range("E30").Value = "  +5.40%  "
$ws.Range("D31").Value = "135.82"
$ws.Range("E31").Value = "  +3.07%  "
$ws.Range("D32").Value = "6.388"
$ws.Range("E32").Value = "  +4.92%  "
$ws.Range("D33").Value = "4.079"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").Value = "13.96"
$ws.Range("E34").Value = "  +7.47%  "
$ws.Range("D35").Value = "0.08818"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").Value = "1.702"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "5.650"
$ws.Range("E37").Value = "  +3.85%  "
$ws.Range("D38").Value = "0.7088"
$ws.Range("E38").Value = "  +13.50%  "
$ws.Range("D39").Value = "9.149"
$ws.Range("E39").Value = "  +6.15%  "
$ws.Range("D40").Value = "0.06523"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.02419"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.2261"
$ws.Range("E42").Value = "  +3.76%  "
$ws.Range("D43").Value = "1.315"
$ws.Range("E43").Value = "  +6.34%  "
$ws.Range("D44").Value = "14.87"
$ws.Range("E44").Value = "  +3.51%  "
$ws.Range("D45").Value = "0.6673"
$ws.Range("E45").Value = "  +10.14%  "
$ws.Range("D46").Value = "0.9990"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "3.960"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").Value = "2.203"
$ws.Range("E48").Value = "  +6.54%  "
$ws.Range("D49").Value = "133.47"
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").Value = "0.07324"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").Value = "81.21"
$ws.Range("E51").Value = "  +4.64%  "
